$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.730.83'
$ws.Range("E2").Value = '  +3.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.443.65'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.19'
$ws.Range("E5").Value = '  +2.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.88'
$ws.Range("E6").Value = '  +4.03%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("E9").Value = '  +1.62%  '
$ws.Range("E10").Value = '  +3.88%  '
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("E13").Value = '  +3.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.87'
$ws.Range("E14").Value = '  +10.76%  '
$ws.Range("E15").Value = '  +6.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.887.63'
$ws.Range("E16").Value = '  +2.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.607.14'
$ws.Range("E17").Value = '  +3.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.448.36'
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.83'
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.96'
$ws.Range("E20").Value = '  +3.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.74'
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.05'
$ws.Range("E23").Value = '  +11.50%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '0.0₆0654'
$ws.Range("E25").Value = '  +126.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.62'
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '624.14'
$ws.Range("E27").Value = '  +9.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.12'
$ws.Range("E28").Value = '  +11.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.47'
$ws.Range("E29").Value = '  +4.40%  '
$ws.Range("D30").Value = '0.0₃0993'
$ws.Range("E30").Value = '  +5.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.563.42'
$ws.Range("E31").Value = '  +2.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.28'
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("E33").Value = '  +6.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.141'
$ws.Range("E34").Value = '  +7.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.86'
$ws.Range("E35").Value = '  +2.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.49'
$ws.Range("E36").Value = '  +2.59%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  +3.02%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.49'
$ws.Range("E39").Value = '  +7.20%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.374'
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '152.90'
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.78'
$ws.Range("E43").Value = '  +11.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.76'
$ws.Range("E44").Value = '  +5.03%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.14'
$ws.Range("E46").Value = '  +29.12%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '144.68'
$ws.Range("E47").Value = '  +2.09%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.61'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.56'
$ws.Range("E49").Value = '  +6.12%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.603'
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("B51").Value = 'Hedera'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0518'
$ws.Range("E51").Value = '  +2.26%  '
